$wb = $excel.ActiveWorkbook

# --- Highlight the internal {R-T;...} template marker cells in red on Sheet3 ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1:A3").Font.Color = 255

# --- Highlight the same marker cells on Sheet4 ---
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("A1:A4").Font.Color = 255

# --- Create a copy of Sheet3 (now red-highlighted), placed right after Sheet3 ---
$ws3.Copy([System.Reflection.Missing]::Value, $ws3)
$newSheet = $wb.Worksheets.Item("Sheet3 (2)")
$newSheet.Name = "Sheet3 ;merge=G0"

# Update the merge markers on the copied sheet from "merge=X" to "merge=G0"
$newSheet.Range("B2").Value = "{R-T-CITYFROM;merge=G0}"
$newSheet.Range("C2").Value = "{R-T-CITYTO;merge=G0}"
$newSheet.Range("D2").Value = "{R-T-CARRNAME;merge=G0}"
$newSheet.Range("E2").Value = "{R-T-CONNID;merge=G0}"

# Add an explanatory note about the new group-based merge syntax
$newSheet.Range("L1").Value = "Make groups like 'G*' to create order based cell merging"
$newSheet.Range("L1").VerticalAlignment = -4108

# Leave the print-preview selection where the user clicked last
$newSheet.Range("T1").Select() | Out-Null

# --- Switch the active tab back to Sheet3 and move the selection ---
$ws3.Activate() | Out-Null
$ws3.Range("M1").Select() | Out-Null
